# Scheduled-runner style update of market/profit figures (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 38.81818
$ws.Range("I11").Value = 38.81818
$ws.Range("K11").Value = 38.81818
$ws.Range("M11").Value = 101.18182
$ws.Range("H34").Value = 7300
$ws.Range("I34").Value = 1625
$ws.Range("J34").Value = 30000
$ws.Range("K34").Value = 1625
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = -1422
$ws.Range("N34").Value = -30406
$ws.Range("H36").Value = 7300
$ws.Range("I36").Value = 1625
$ws.Range("J36").Value = 30000
$ws.Range("K36").Value = 1625
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = -910
$ws.Range("N36").Value = -31430
$ws.Range("H40").Value = 2857.1428
$ws.Range("I40").Value = 3160
$ws.Range("J40").Value = 2100
$ws.Range("K40").Value = 3160
$ws.Range("L40").Value = 2100
$ws.Range("M40").Value = -2985
$ws.Range("N40").Value = -2450
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30924
$ws.Range("H47").Value = 13643.5
$ws.Range("I47").Value = 10500
$ws.Range("J47").Value = 16787
$ws.Range("K47").Value = 10500
$ws.Range("L47").Value = 16787
$ws.Range("M47").Value = -9528
$ws.Range("N47").Value = -18731
$ws.Range("H76").Value = 3495.5386
$ws.Range("I76").Value = 2833.3333
$ws.Range("K76").Value = 2833.3333
$ws.Range("M76").Value = -2518.3333
$ws.Range("H79").Value = 3495.5386
$ws.Range("I79").Value = 2833.3333
$ws.Range("K79").Value = 2833.3333
$ws.Range("M79").Value = -1741.3333
$ws.Range("H112").Value = 1439.4889
$ws.Range("I112").Value = 598
$ws.Range("K112").Value = 1794
$ws.Range("M112").Value = -686
$ws.Range("H131").Value = 70050.734
$ws.Range("I131").Value = 113817.336
$ws.Range("J131").Value = 4400.8335
$ws.Range("K131").Value = 341452.008
$ws.Range("L131").Value = 13202.5005
$ws.Range("M131").Value = -336412.008
$ws.Range("N131").Value = -23282.5005
$ws.Range("H138").Value = 4100.956
$ws.Range("J138").Value = 4362.271
$ws.Range("L138").Value = 13086.813
$ws.Range("N138").Value = -23366.813
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26343232
$ws.Range("I32").Value = 55569160
$ws.Range("J32").Value = 39900
$ws.Range("K32").Value = 55569160
$ws.Range("L32").Value = 39900
$ws.Range("M32").Value = -55568873
$ws.Range("N32").Value = -40474
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H39").Value = 2658
$ws.Range("I39").Value = 2658
$ws.Range("K39").Value = 2658
$ws.Range("M39").Value = -2138
$ws.Range("H56").Value = 10500
$ws.Range("J56").Value = 10500
$ws.Range("L56").Value = 10500
$ws.Range("N56").Value = -11984
$ws.Range("H108").Value = 30750
$ws.Range("J108").Value = 30750
$ws.Range("L108").Value = 30750
$ws.Range("N108").Value = -38430
$ws.Range("H122").Value = 1535.7273
$ws.Range("I122").Value = 1389.3
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4167.9
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1717.9
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 1413.0638
$ws.Range("I132").Value = 864.6129
$ws.Range("J132").Value = 2475.6875
$ws.Range("K132").Value = 2593.8387
$ws.Range("L132").Value = 7427.0625
$ws.Range("M132").Value = -63.83869999999979
$ws.Range("N132").Value = -12487.0625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2065.7083
$ws.Range("I86").Value = 2044.8
$ws.Range("J86").Value = 2100.5557
$ws.Range("K86").Value = 2044.8
$ws.Range("L86").Value = 2100.5557
$ws.Range("M86").Value = -921.8
$ws.Range("N86").Value = -4346.5557
$ws.Range("H89").Value = 2065.7083
$ws.Range("I89").Value = 2044.8
$ws.Range("J89").Value = 2100.5557
$ws.Range("K89").Value = 10224
$ws.Range("L89").Value = 10502.7785
$ws.Range("M89").Value = -4608
$ws.Range("N89").Value = -21734.7785
$ws.Range("H132").Value = 40597.5
$ws.Range("J132").Value = 40597.5
$ws.Range("L132").Value = 40597.5
$ws.Range("N132").Value = -50717.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 17350
$ws.Range("J109").Value = 17350
$ws.Range("L109").Value = 17350
$ws.Range("N109").Value = -19430
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 290
$ws.Range("I86").Value = 313.33334
$ws.Range("J86").Value = 220
$ws.Range("K86").Value = 940.0000200000001
$ws.Range("L86").Value = 660
$ws.Range("M86").Value = 245.9999799999999
$ws.Range("N86").Value = -3032
$ws.Range("H89").Value = 290
$ws.Range("I89").Value = 313.33334
$ws.Range("J89").Value = 220
$ws.Range("K89").Value = 2820.00006
$ws.Range("L89").Value = 1980
$ws.Range("M89").Value = 3107.99994
$ws.Range("N89").Value = -13836
$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -252
$ws.Range("N92").ClearContents()
$ws.Range("H107").Value = 625269.4399999999
$ws.Range("I107").Value = 244.88889
$ws.Range("J107").Value = 1428872.4
$ws.Range("K107").Value = 734.6666700000001
$ws.Range("L107").Value = 4286617.199999999
$ws.Range("M107").Value = 1185.33333
$ws.Range("N107").Value = -4290457.199999999
$ws.Range("H132").Value = 3929334
$ws.Range("I132").Value = 848.625
$ws.Range("J132").Value = 6096774.5
$ws.Range("K132").Value = 7637.625
$ws.Range("L132").Value = 54870970.5
$ws.Range("M132").Value = -5107.625
$ws.Range("N132").Value = -54876030.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 3000
$ws.Range("I52").Value = 3000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -2741
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 7742.857
$ws.Range("I70").Value = 9240
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 9240
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -8970
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 7742.857
$ws.Range("I73").Value = 9240
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 9240
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -8304
$ws.Range("N73").Value = -5872
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -13900
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6111.619
$ws.Range("I136").Value = 8067.7856
$ws.Range("J136").Value = 2199.2856
$ws.Range("K136").Value = 24203.3568
$ws.Range("L136").Value = 6597.8568
$ws.Range("M136").Value = -21653.3568
$ws.Range("N136").Value = -11697.8568
